$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (State) to make room for "Process"
$ws.Columns.Item(7).Insert()

# Header/data for new column - copy formatting from neighbor (H1:H3) then set values
$ws.Range("H1:H3").Copy()
$ws.Range("G1:G3").PasteSpecial(-4122)

$ws.Cells.Item(1, 7).Value = "Process"
$ws.Cells.Item(2, 7).Value = "Search"
$ws.Cells.Item(3, 7).Value = "Typing"

# Update OrderID values
$ws.Cells.Item(2, 2).Value = 1213286
$ws.Cells.Item(3, 2).Value = 2193289

# Column widths for the new column and following ones
$ws.Columns.Item(7).ColumnWidth = 15.2

# Update selection
$ws.Range("G10").Select()
